$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "TempR"
$ws.Range("F2").Value = "HumR"

$ws.Range("B2").Copy()
$ws.Range("E2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E2:F2").Borders.LineStyle = 0

# Insert a row above row 17, forcing row 17 data to shift to row 18, row 18->19
$ws.Rows("17:17").Insert()
# Now delete the blank inserted row, shifting back up
$ws.Rows("17:17").Delete()

$ws.Range("E3").Select()
